# Loop over excel file and get email data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the mailto: hyperlink that used to live on A1.
$ws.Hyperlinks.Delete()

# Replace the header with a plain "Email" label.
$ws.Range("A1").Value = "Email"

# Row no longer needs the tall wrapped-text height - reset it to the sheet default.
$ws.Rows(1).AutoFit()

# Add the looped-over test rows below the header (write A3 first so the
# shared-string table interns "test 2" before "test 1").
$ws.Range("A3").Value = "test 2"
$ws.Range("A2").Value = "test 1"

# Move the selection to B1 to match the post-edit view state.
$ws.Range("B1").Select()
